$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'243.69"
$ws.Range("D4").Value = "'5.388"
$ws.Range("D5").Value = "'0.05958"
$ws.Range("D7").Value = "'6.498"
$ws.Range("D8").Value = "'0.8108"
$ws.Range("D9").Value = "'0.9296"
$ws.Range("D10").Value = "'0.1435"
$ws.Range("D11").Value = "'0.07417"
$ws.Range("D12").Value = "'0.03270"
$ws.Range("D13").Value = "'0.03080"
$ws.Range("D14").Value = "'0.09363"
$ws.Range("D15").Value = "'3.844"
$ws.Range("D16").Value = "'0.001579"
$ws.Range("D17").Value = "'0.04704"
$ws.Range("D18").Value = "'0.0005988"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.005939"
$ws.Range("D20").Value = "'0.001261"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("D21").Value = "'0.004793"
$ws.Range("D23").Value = "'3.575"
$ws.Range("D25").Value = "'0.3240"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"
$ws.Range("D40").Value = "'0.03936"
$ws.Range("D41").Value = "'0.006338"
$ws.Range("D42").Value = "'0.1077"
$ws.Range("D43").Value = "'0.003499"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("D44").Value = "'0.008951"
$ws.Range("D45").Value = "'0.00005179"
$ws.Range("D47").Value = "'0.6848"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.0001999"
